$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01135533333333333
$ws.Range("H2").Value = 0.034066
$ws.Range("M2").Value = 68.46613766666667
$ws.Range("N2").Value = 205.398413
$ws.Range("O2").Value = 0.4719163120948675
$ws.Range("P2").Value = 0.4719163120948675
$ws.Range("Q2").Value = 0.7774558152508889
$ws.Range("R2").Value = 6.997102337258
$ws.Range("S2").Value = 0.4719163120948675
$ws.Range("T2").Value = 0.4719163120948675

# Row 3
$ws.Range("G3").Value = 0.01135533333333333
$ws.Range("H3").Value = 0.034066
$ws.Range("M3").Value = 9.278736333333333
$ws.Range("O3").Value = 0.06395551407683932
$ws.Range("P3").Value = 0.06395551407683933
$ws.Range("Q3").Value = 0.1053631439771111
$ws.Range("R3").Value = 0.948268295794
$ws.Range("S3").Value = 0.06395551407683932
$ws.Range("T3").Value = 0.06395551407683933

# Row 4
$ws.Range("G4").Value = 0.01135533333333333
$ws.Range("H4").Value = 0.034066
$ws.Range("M4").Value = 67.336226
$ws.Range("N4").Value = 202.008678
$ws.Range("O4").Value = 0.4641281738282933
$ws.Range("P4").Value = 0.4641281738282933
$ws.Range("Q4").Value = 0.7646252916386667
$ws.Range("R4").Value = 6.881627624748
$ws.Range("S4").Value = 0.4641281738282933
$ws.Range("T4").Value = 0.4641281738282933

$wb.Save()
